{"js": "// Search for the run containing fldChar begin, get its range to end of paragraph\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst p = paras.items[1];\nconst pRange = p.getRange();\n// Use search to find \"Artifact1\" text which is after the field begin\nconst results = pRange.search(\"Artifact1\", {matchCase:true});\nresults.load(\"items\");\nawait context.sync();\nreturn JSON.stringify(results.items.length);\n", "ps1": "$d = $word.ActiveDocument\nWrite-Output $d.CurrentRsid\n"}
